$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace paragraph 1 ("NOMOR ANTRIAN") with an edited copy of itself
#    (border w:space 1 -> 0, drop explicit sz/szCs=18 overrides) followed by
#    a brand-new paragraph ("Menuju Nurse Station") that carries the same
#    border + sz=16/szCs=16.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="33F34423" w14:textId="6E314D08" w:rsidR="00980477" w:rsidRPr="00980477" w:rsidRDefault="00B51030" w:rsidP="00F97C4C" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:pPr>
<w:pBdr><w:bottom w:val="double" w:sz="6" w:space="0" w:color="auto"/></w:pBdr>
<w:spacing w:line="276" w:lineRule="auto"/>
<w:jc w:val="center"/>
<w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia"/><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr>
</w:pPr>
<w:r w:rsidRPr="00F4028F"><w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia"/><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">NOMOR </w:t></w:r>
<w:r w:rsidR="00F4028F" w:rsidRPr="00F4028F"><w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia"/><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>ANTRIAN</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pBdr><w:bottom w:val="double" w:sz="6" w:space="0" w:color="auto"/></w:pBdr>
<w:spacing w:line="276" w:lineRule="auto"/>
<w:jc w:val="center"/>
<w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>Menuju</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Nurse Station</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Add <w10:wrap anchorx="margin" anchory="margin"/> into every watermark
#    VML shape (one per header part) right after its <v:imagedata .../>.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("gain=`"19661f`" blacklevel=`"22938f`"/></v:shape>", $false, $false, $false, $false, $false, $true, 1, $false, "gain=`"19661f`" blacklevel=`"22938f`"/><w10:wrap anchorx=`"margin`" anchory=`"margin`"/></v:shape>", 2)
